$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 319245
$ws.Cells.Item(2, 4).Value = 406871348
$ws.Cells.Item(4, 3).Value = 315
$ws.Cells.Item(4, 4).Value = 450707
$ws.Cells.Item(8, 3).Value = 858
$ws.Cells.Item(8, 4).Value = 1262408
$ws.Cells.Item(10, 3).Value = 116637
$ws.Cells.Item(10, 4).Value = 170910039
$ws.Cells.Item(12, 3).Value = 59099
$ws.Cells.Item(12, 4).Value = 85293945
$ws.Cells.Item(16, 3).Value = 4000
$ws.Cells.Item(16, 4).Value = 5676575
$ws.Cells.Item(20, 3).Value = 6589
$ws.Cells.Item(20, 4).Value = 9192134
$ws.Cells.Item(22, 3).Value = 77128
$ws.Cells.Item(22, 4).Value = 96191150
$ws.Cells.Item(28, 3).Value = 32352
$ws.Cells.Item(28, 4).Value = 47359432
$ws.Cells.Item(30, 3).Value = 11432
$ws.Cells.Item(30, 4).Value = 16441040
$ws.Cells.Item(33, 3).Value = 1559
$ws.Cells.Item(33, 4).Value = 2189807
$ws.Cells.Item(35, 3).Value = 1809
$ws.Cells.Item(35, 4).Value = 2553671
$ws.Cells.Item(36, 3).Value = 96746
$ws.Cells.Item(36, 4).Value = 121789127
$ws.Cells.Item(42, 3).Value = 901
$ws.Cells.Item(42, 4).Value = 1326185
$ws.Cells.Item(44, 3).Value = 44251
$ws.Cells.Item(44, 4).Value = 64849724
$ws.Cells.Item(46, 3).Value = 9098
$ws.Cells.Item(46, 4).Value = 13054796
$ws.Cells.Item(48, 3).Value = 1401
$ws.Cells.Item(48, 4).Value = 1945010
$ws.Cells.Item(51, 3).Value = 2285
$ws.Cells.Item(51, 4).Value = 3188628
$ws.Cells.Item(52, 3).Value = 68780
$ws.Cells.Item(52, 4).Value = 86282932
$ws.Cells.Item(54, 3).Value = 44
$ws.Cells.Item(54, 4).Value = 60697
$ws.Cells.Item(59, 3).Value = 28084
$ws.Cells.Item(59, 4).Value = 41189366
$ws.Cells.Item(62, 3).Value = 11051
$ws.Cells.Item(62, 4).Value = 15979479
$ws.Cells.Item(68, 3).Value = 1457
$ws.Cells.Item(68, 4).Value = 2041565
$ws.Cells.Item(70, 3).Value = 20414
$ws.Cells.Item(70, 4).Value = 26742300
$ws.Cells.Item(73, 3).Value = 59
$ws.Cells.Item(73, 4).Value = 86073
$ws.Cells.Item(74, 3).Value = 7560
$ws.Cells.Item(74, 4).Value = 11069064
$ws.Cells.Item(76, 3).Value = 5104
$ws.Cells.Item(76, 4).Value = 7409384
$ws.Cells.Item(78, 3).Value = 271
$ws.Cells.Item(78, 4).Value = 380173
$ws.Cells.Item(79, 3).Value = 140259
$ws.Cells.Item(79, 4).Value = 174919792
$ws.Cells.Item(81, 3).Value = 88
$ws.Cells.Item(81, 4).Value = 123384
$ws.Cells.Item(83, 3).Value = 425
$ws.Cells.Item(83, 4).Value = 620771
$ws.Cells.Item(85, 3).Value = 63333
$ws.Cells.Item(85, 4).Value = 92824331
$ws.Cells.Item(88, 3).Value = 29592
$ws.Cells.Item(88, 4).Value = 42808751
$ws.Cells.Item(91, 3).Value = 2803
$ws.Cells.Item(91, 4).Value = 3959864
$ws.Cells.Item(92, 3).Value = 32726
$ws.Cells.Item(92, 4).Value = 44337200
$ws.Cells.Item(95, 3).Value = 27
$ws.Cells.Item(95, 4).Value = 38814
$ws.Cells.Item(96, 3).Value = 7903
$ws.Cells.Item(96, 4).Value = 11619686
$ws.Cells.Item(98, 3).Value = 7243
$ws.Cells.Item(98, 4).Value = 10503011
$ws.Cells.Item(100, 3).Value = 530
$ws.Cells.Item(100, 4).Value = 753216
$ws.Cells.Item(102, 3).Value = 9207
$ws.Cells.Item(102, 4).Value = 13033754
$ws.Cells.Item(104, 3).Value = 2338
$ws.Cells.Item(104, 4).Value = 3545597
$ws.Cells.Item(106, 3).Value = 3113
$ws.Cells.Item(106, 4).Value = 4660982
$ws.Cells.Item(108, 3).Value = 133
$ws.Cells.Item(108, 4).Value = 195120
$ws.Cells.Item(109, 3).Value = 179
$ws.Cells.Item(109, 4).Value = 253843
$ws.Cells.Item(110, 3).Value = 140868
$ws.Cells.Item(110, 4).Value = 174213229
$ws.Cells.Item(116, 3).Value = 52526
$ws.Cells.Item(116, 4).Value = 76996340
$ws.Cells.Item(118, 3).Value = 26881
$ws.Cells.Item(118, 4).Value = 38943412
$ws.Cells.Item(119, 3).Value = 1309
$ws.Cells.Item(119, 4).Value = 1790634
$ws.Cells.Item(122, 3).Value = 2240
$ws.Cells.Item(122, 4).Value = 3144983
$ws.Cells.Item(124, 3).Value = 504630
$ws.Cells.Item(124, 4).Value = 665993506
$ws.Cells.Item(126, 3).Value = 212
$ws.Cells.Item(126, 4).Value = 312009
$ws.Cells.Item(129, 3).Value = 1371
$ws.Cells.Item(129, 4).Value = 2032811
$ws.Cells.Item(130, 3).Value = 32
$ws.Cells.Item(130, 4).Value = 42010
$ws.Cells.Item(131, 3).Value = 206643
$ws.Cells.Item(131, 4).Value = 303763413
$ws.Cells.Item(132, 3).Value = 393
$ws.Cells.Item(132, 4).Value = 586250
$ws.Cells.Item(134, 3).Value = 180157
$ws.Cells.Item(134, 4).Value = 261897201
$ws.Cells.Item(137, 3).Value = 2846
$ws.Cells.Item(137, 4).Value = 3998456
$ws.Cells.Item(139, 3).Value = 6295
$ws.Cells.Item(139, 4).Value = 8894620
$ws.Cells.Item(142, 3).Value = 44349
$ws.Cells.Item(142, 4).Value = 59209720
$ws.Cells.Item(148, 3).Value = 13996
$ws.Cells.Item(148, 4).Value = 20525287
$ws.Cells.Item(149, 3).Value = 3735
$ws.Cells.Item(149, 4).Value = 5386097
$ws.Cells.Item(152, 3).Value = 400
$ws.Cells.Item(152, 4).Value = 575216
$ws.Cells.Item(154, 3).Value = 378
$ws.Cells.Item(154, 4).Value = 533751
$ws.Cells.Item(155, 3).Value = 17472
$ws.Cells.Item(155, 4).Value = 23088912
$ws.Cells.Item(159, 3).Value = 7118
$ws.Cells.Item(159, 4).Value = 10352086
$ws.Cells.Item(161, 3).Value = 4960
$ws.Cells.Item(161, 4).Value = 7137333
$ws.Cells.Item(163, 3).Value = 276
$ws.Cells.Item(163, 4).Value = 381431
$ws.Cells.Item(166, 3).Value = 16323
$ws.Cells.Item(166, 4).Value = 24422175
$ws.Cells.Item(167, 3).Value = 1839
$ws.Cells.Item(167, 4).Value = 2812613
$ws.Cells.Item(168, 3).Value = 243
$ws.Cells.Item(168, 4).Value = 366802
$ws.Cells.Item(170, 3).Value = 55
$ws.Cells.Item(170, 4).Value = 85190
$ws.Cells.Item(171, 3).Value = 95
$ws.Cells.Item(171, 4).Value = 154449
$ws.Cells.Item(172, 3).Value = 87367
$ws.Cells.Item(172, 4).Value = 109273205
$ws.Cells.Item(176, 3).Value = 14
$ws.Cells.Item(176, 4).Value = 20820
$ws.Cells.Item(179, 3).Value = 33738
$ws.Cells.Item(179, 4).Value = 49476963
$ws.Cells.Item(181, 3).Value = 12953
$ws.Cells.Item(181, 4).Value = 18715230
$ws.Cells.Item(183, 3).Value = 1244
$ws.Cells.Item(183, 4).Value = 1741896
$ws.Cells.Item(185, 3).Value = 1633
$ws.Cells.Item(185, 4).Value = 2294635
$ws.Cells.Item(187, 3).Value = 237683
$ws.Cells.Item(187, 4).Value = 295485117
$ws.Cells.Item(189, 3).Value = 168
$ws.Cells.Item(189, 4).Value = 242236
$ws.Cells.Item(193, 3).Value = 873
$ws.Cells.Item(193, 4).Value = 1283845
$ws.Cells.Item(195, 3).Value = 86314
$ws.Cells.Item(195, 4).Value = 126525712
$ws.Cells.Item(198, 3).Value = 32893
$ws.Cells.Item(198, 4).Value = 47343968
$ws.Cells.Item(201, 3).Value = 5101
$ws.Cells.Item(201, 4).Value = 7261436
$ws.Cells.Item(204, 3).Value = 4835
$ws.Cells.Item(204, 4).Value = 6696033
$ws.Cells.Item(207, 3).Value = 262816
$ws.Cells.Item(207, 4).Value = 325294041
$ws.Cells.Item(214, 3).Value = 612
$ws.Cells.Item(214, 4).Value = 890878
$ws.Cells.Item(216, 3).Value = 94829
$ws.Cells.Item(216, 4).Value = 138730121
$ws.Cells.Item(217, 3).Value = 88
$ws.Cells.Item(217, 4).Value = 131199
$ws.Cells.Item(219, 3).Value = 51221
$ws.Cells.Item(219, 4).Value = 74033583
$ws.Cells.Item(222, 3).Value = 4667
$ws.Cells.Item(222, 4).Value = 6553355
$ws.Cells.Item(225, 3).Value = 5698
$ws.Cells.Item(225, 4).Value = 7880764
$ws.Cells.Item(228, 3).Value = 105830
$ws.Cells.Item(228, 4).Value = 132390997
$ws.Cells.Item(231, 3).Value = 10
$ws.Cells.Item(231, 4).Value = 12647
$ws.Cells.Item(233, 3).Value = 564
$ws.Cells.Item(233, 4).Value = 823939
$ws.Cells.Item(235, 3).Value = 49299
$ws.Cells.Item(235, 4).Value = 72217501
$ws.Cells.Item(237, 3).Value = 12340
$ws.Cells.Item(237, 4).Value = 17745369
$ws.Cells.Item(239, 3).Value = 1890
$ws.Cells.Item(239, 4).Value = 2708882
$ws.Cells.Item(241, 3).Value = 2492
$ws.Cells.Item(241, 4).Value = 3484565
$ws.Cells.Item(242, 3).Value = 256350
$ws.Cells.Item(242, 4).Value = 323705322
$ws.Cells.Item(244, 3).Value = 249
$ws.Cells.Item(244, 4).Value = 357957
$ws.Cells.Item(248, 3).Value = 825
$ws.Cells.Item(248, 4).Value = 1211563
$ws.Cells.Item(249, 3).Value = 10
$ws.Cells.Item(249, 4).Value = 15000
$ws.Cells.Item(250, 3).Value = 95424
$ws.Cells.Item(250, 4).Value = 139818032
$ws.Cells.Item(253, 3).Value = 64672
$ws.Cells.Item(253, 4).Value = 93718520
$ws.Cells.Item(255, 3).Value = 2408
$ws.Cells.Item(255, 4).Value = 3396674
$ws.Cells.Item(258, 3).Value = 4566
$ws.Cells.Item(258, 4).Value = 6411613
